$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text so numeric-looking strings
# (e.g. "233.48", "14.80") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.810.38"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.076.97"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "233.48"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "59.11"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "14.80"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "21.19"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "2.118.58"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "37.705.86"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "71.58"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "0.0₃0842"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "228.30"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").Value = "171.02"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").Value = "9.17"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  -3.34%  "
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "4.76"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").Value = "0.0634"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "2.48"
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("D35").Value = "1.83"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "0.0977"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").Value = "99.41"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").Value = "16.67"
$ws.Range("E43").Value = "  +5.65%  "
$ws.Range("D44").Value = "1.440.35"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").Value = "7.41"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").Value = "2.267.87"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  +0.46%  "

# Remove the temporary text-format styling so cell styles are left unchanged,
# while keeping the values stored as text strings.
$ws.Range("D2:D51").ClearFormats()
